$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42 (pushes the existing rows 42-59 down to 43-60)
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with the new weekly price entry
$ws.Cells.Item(42, 1).Value2  = 10
$ws.Cells.Item(42, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(42, 3).Value2  = "La Araucanía"
$ws.Cells.Item(42, 4).Value2  = 44518
$ws.Cells.Item(42, 5).Value2  = 9
$ws.Cells.Item(42, 6).Value2  = "Fruta"
$ws.Cells.Item(42, 7).Value2  = 100101
$ws.Cells.Item(42, 8).Value2  = "Berries"
$ws.Cells.Item(42, 9).Value2  = 100101001
$ws.Cells.Item(42, 10).Value2 = "Arándano (blue)"
$ws.Cells.Item(42, 11).Value2 = "Sin especificar"
$ws.Cells.Item(42, 12).Value2 = "Primera"
$ws.Cells.Item(42, 13).Value2 = 80
$ws.Cells.Item(42, 14).Value2 = 3500
$ws.Cells.Item(42, 15).Value2 = 3500
$ws.Cells.Item(42, 16).Value2 = 3500
$ws.Cells.Item(42, 17).Value2 = "$/kilo"
$ws.Cells.Item(42, 18).Value2 = "Región del Maule"
$ws.Cells.Item(42, 19).Value2 = 3500
$ws.Cells.Item(42, 20).Value2 = 1
